$d = $word.ActiveDocument

# Locate the first of the three reference paragraphs that need to be
# replaced (the one containing the gist.github.com hyperlink).
$targetText = "https://gist.github.com/ftvs/5822103"
$firstParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*$targetText*") {
        $firstParaIndex = $i
        break
    }
}

$p = $d.Paragraphs.Item($firstParaIndex)
$r = $p.Range
$r.Collapse(1)

$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:proofErr w:type="spellStart"/>
<w:r><w:t>ftvs</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> (2014) </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:i/></w:rPr><w:t>CameraShake.cs</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">. </w:t></w:r>
<w:r><w:t xml:space="preserve">Available at: </w:t></w:r>
<w:r><w:t>https://gist.github.com/ftvs/5822103</w:t></w:r>
<w:r><w:t xml:space="preserve"> (Accessed: 8 December 2019).</w:t></w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:proofErr w:type="spellStart"/>
<w:r><w:t>Kakkuonhyvaa</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> (2016) </w:t></w:r>
<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">Scrolling sprite texture. </w:t></w:r>
<w:r><w:t xml:space="preserve">Available at: </w:t></w:r>
<w:r><w:t>https://answers.unity.com/questions/1080218/scrolling-sprite-texture.html</w:t></w:r>
<w:r><w:t xml:space="preserve"> (Accessed: 9 December 2019).</w:t></w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:r><w:t xml:space="preserve">Resistance Studio (2018) </w:t></w:r>
<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">Making your Pixel Art Game </w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:rPr><w:i/></w:rPr><w:t>look</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> Pixel Perfect in Unity3D. </w:t></w:r>
<w:r><w:t xml:space="preserve">Available at: </w:t></w:r>
<w:r><w:t>https://hackernoon.com/making-your-pixel-art-game-look-pixel-perfect-in-unity3d-3534963cad1d</w:t></w:r>
<w:r><w:t xml:space="preserve"> (Accessed: 19 December 2019).</w:t></w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:proofErr w:type="spellStart"/>
<w:r><w:lastRenderedPageBreak/><w:t>Brackeys</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> (2019) </w:t></w:r>
<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">CONTROLLER INPUT in Unity. </w:t></w:r>
<w:r><w:t xml:space="preserve">14 July. </w:t></w:r>
<w:r><w:t xml:space="preserve">Available at: </w:t></w:r>
<w:r><w:t>https://www.youtube.com/watch?v=p-3S73MaDP8</w:t></w:r>
<w:r><w:t xml:space="preserve"> (Accessed: 24 December 2019).</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
'@

$r.InsertXML($newXml)

# Remove the now-orphaned two paragraphs that used to hold the
# answers.unity.com and hackernoon.com hyperlinks (they were pushed
# further down the document by the insert above, and are no longer
# needed since their content was folded into the new paragraphs).
$answersText = "https://answers.unity.com/questions/1080218/scrolling-sprite-texture.html"
$hackernoonText = "https://hackernoon.com/making-your-pixel-art-game-look-pixel-perfect-in-unity3d-3534963cad1d"

$answersIdx = -1
$hackernoonIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text.Trim()
    if ($answersIdx -eq -1 -and $txt -eq $answersText) {
        $answersIdx = $i
    }
    if ($hackernoonIdx -eq -1 -and $txt -eq $hackernoonText) {
        $hackernoonIdx = $i
    }
}

$pAnswers = $d.Paragraphs.Item($answersIdx)
$pHackernoon = $d.Paragraphs.Item($hackernoonIdx)
$delRange = $d.Range($pAnswers.Range.Start, $pHackernoon.Range.End)
$delRange.Delete()
